# updated script and files in csv instead of xlsx
# - row2 (modelPath) and row4 (modelPath) point to a new source path
# - the two affected records' identifier UUIDs are bumped to a new version
# - a duplicate of the (updated) row4 record is appended as row5, but keeping
#   the *old* modelPath value in column A (the source file for that record
#   hadn't moved yet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture row 4's current ("before edit") values across all used columns ---
# (row4 is going to be edited in place; we need its pre-edit contents to build
#  the new row5 further down)
$lastCol = $ws.UsedRange.Columns.Count
$oldRow4 = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $oldRow4 += ,$ws.Cells.Item(4, $c).Value2
}

# --- row 2 edits: new modelPath + bumped identifier uuid ---
$ws.Range("A2").Value2 = "/home/libotadmin/NewYork"
$ws.Range("C2").Value2 = "b6645aa5-4134-50f3-8cbc-faa0518c21bb"

# --- row 3 edits: bumped identifier uuid only ---
$ws.Range("C3").Value2 = "Dc61d24d-7426-5090-a48f-06c13be98b85"

# --- row 4 edits: new modelPath + bumped identifier uuid ---
$ws.Range("A4").Value2 = "/home/libotadmin/NewYork"
$ws.Range("C4").Value2 = "Deda97ad-0912-5024-929c-02beba91c01d"

# --- row 2 & 4 column A formatting: wrap text + taller row ---
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 15

$ws.Range("A4").Font.Name = "Calibri"
$ws.Range("A4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 15

# --- new row 5: same as the (now updated) row 4, except column A keeps the
#     record's previous modelPath value ---
for ($c = 1; $c -le $lastCol; $c++) {
    if ($c -eq 1) {
        $ws.Cells.Item(5, $c).Value2 = $oldRow4[$c - 1]
    } else {
        $ws.Cells.Item(5, $c).Value2 = $ws.Cells.Item(4, $c).Value2
    }
}

# --- selection ends on C5, matching the saved state after the edit ---
$ws.Range("C5").Select()
